$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing amounts
$ws.Range("B3").Value = 20
$ws.Range("B5").Value = 7.5

# Add new row 8: Prints, 38, 10/26/2016 (reuse formatting from row 7)
$ws.Range("A8").Value = "Prints"
$ws.Range("B8").Value = 38
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C8").Value2 = 42669

$excel.CutCopyMode = 0

# Update selection to C9
$ws.Range("C9").Select()
